$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-26 Sunday", "2025-01-27 Monday"),
    @("194÷6=", "219÷9="),
    @("138÷9=", "888÷4="),
    @("492÷4=", "459÷6="),
    @("292÷4=", "809÷3="),
    @("406÷9=", "210÷2="),
    @("457÷4=", "850÷8="),
    @("449÷6=", "716÷5="),
    @("354÷2=", "776÷4="),
    @("868÷8=", "969÷6="),
    @("747÷8=", "237÷8="),
    @("533÷5=", "893÷9="),
    @("266÷5=", "197÷5="),
    @("880÷8=", "132÷4="),
    @("315÷2=", "822÷4="),
    @("434÷3=", "240÷2="),
    @("164÷9=", "702÷9="),
    @("612÷7=", "510÷8="),
    @("272÷3=", "419÷8="),
    @("167÷5=", "781÷9="),
    @("161÷6=", "984÷5="),
    @("366÷6=", "535÷6="),
    @("362÷7=", "188÷5="),
    @("840÷8=", "959÷3="),
    @("566÷4=", "449÷5="),
    @("602÷3=", "882÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
